$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Cells.Item(2, 10).Value = 6989
$ws.Cells.Item(3, 10).Value = 7388
$ws.Cells.Item(4, 3).Value = 1844
$ws.Cells.Item(4, 5).Value = 2015
$ws.Cells.Item(4, 10).Value = 1609
$ws.Cells.Item(5, 10).Value = 578
$ws.Cells.Item(6, 9).Value = 8965
$ws.Cells.Item(6, 10).Value = 9995
$ws.Cells.Item(7, 3).Value = 28388
$ws.Cells.Item(7, 5).Value = 26020
$ws.Cells.Item(7, 9).Value = 26232
$ws.Cells.Item(7, 10).Value = 26559

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Cells.Item(2, 10).Value = 69
$ws.Cells.Item(3, 10).Value = 45
$ws.Cells.Item(6, 10).Value = 262
$ws.Cells.Item(7, 10).Value = 392

$ws = $wb.Worksheets.Item("Austin")
$ws.Cells.Item(2, 10).Value = 444
$ws.Cells.Item(3, 10).Value = 494
$ws.Cells.Item(6, 10).Value = 599
$ws.Cells.Item(7, 10).Value = 1667

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Cells.Item(4, 10).Value = 21
$ws.Cells.Item(6, 10).Value = 142
$ws.Cells.Item(7, 10).Value = 532

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Cells.Item(3, 10).Value = 136
$ws.Cells.Item(4, 10).Value = 18

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Cells.Item(3, 10).Value = 279
$ws.Cells.Item(4, 10).Value = 32
$ws.Cells.Item(6, 9).Value = 243
$ws.Cells.Item(7, 9).Value = 809
$ws.Cells.Item(7, 10).Value = 823

$ws = $wb.Worksheets.Item("New City")
$ws.Cells.Item(2, 10).Value = 191
$ws.Cells.Item(3, 10).Value = 186
$ws.Cells.Item(5, 10).Value = 18
$ws.Cells.Item(7, 10).Value = 671

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Cells.Item(6, 10).Value = 103
$ws.Cells.Item(7, 10).Value = 404

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Cells.Item(3, 10).Value = 36
$ws.Cells.Item(7, 10).Value = 93

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Cells.Item(5, 10).Value = 80
$ws.Cells.Item(6, 10).Value = 200
$ws.Cells.Item(8, 10).Value = 1667
$ws.Cells.Item(11, 10).Value = 468
$ws.Cells.Item(14, 10).Value = 143
$ws.Cells.Item(15, 10).Value = 325
$ws.Cells.Item(19, 10).Value = 770
$ws.Cells.Item(20, 10).Value = 562
$ws.Cells.Item(22, 10).Value = 63
$ws.Cells.Item(23, 10).Value = 243
$ws.Cells.Item(29, 10).Value = 1426
$ws.Cells.Item(30, 10).Value = 93
$ws.Cells.Item(34, 10).Value = 120
$ws.Cells.Item(36, 10).Value = 361
$ws.Cells.Item(37, 9).Value = 809
$ws.Cells.Item(37, 10).Value = 823
$ws.Cells.Item(39, 10).Value = 19
$ws.Cells.Item(42, 10).Value = 1147
$ws.Cells.Item(43, 10).Value = 227
$ws.Cells.Item(46, 10).Value = 90
$ws.Cells.Item(51, 10).Value = 326
$ws.Cells.Item(53, 10).Value = 392
$ws.Cells.Item(54, 10).Value = 514
$ws.Cells.Item(55, 10).Value = 419
$ws.Cells.Item(58, 10).Value = 15
$ws.Cells.Item(60, 10).Value = 156
$ws.Cells.Item(63, 3).Value = 273
$ws.Cells.Item(63, 10).Value = 86
$ws.Cells.Item(64, 10).Value = 175
$ws.Cells.Item(65, 10).Value = 671
$ws.Cells.Item(67, 5).Value = 1129
$ws.Cells.Item(67, 10).Value = 986
$ws.Cells.Item(71, 10).Value = 86
$ws.Cells.Item(80, 10).Value = 45
$ws.Cells.Item(83, 10).Value = 532
$ws.Cells.Item(84, 10).Value = 221
$ws.Cells.Item(86, 10).Value = 168
$ws.Cells.Item(87, 10).Value = 86
$ws.Cells.Item(89, 10).Value = 333
$ws.Cells.Item(90, 10).Value = 279
$ws.Cells.Item(91, 10).Value = 306
$ws.Cells.Item(94, 10).Value = 290
$ws.Cells.Item(96, 10).Value = 288
$ws.Cells.Item(99, 10).Value = 404
$ws.Cells.Item(101, 3).Value = 28388
$ws.Cells.Item(101, 5).Value = 26020
$ws.Cells.Item(101, 9).Value = 26232
$ws.Cells.Item(101, 10).Value = 26559

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Cells.Item(2, 10).Value = 253
$ws.Cells.Item(3, 10).Value = 367
$ws.Cells.Item(4, 5).Value = 56
$ws.Cells.Item(7, 5).Value = 1129
$ws.Cells.Item(7, 10).Value = 986

$ws = $wb.Worksheets.Item("South Deering")
$ws.Cells.Item(6, 10).Value = 71
$ws.Cells.Item(7, 10).Value = 221

$ws = $wb.Worksheets.Item("Loop")
$ws.Cells.Item(4, 10).Value = 38
$ws.Cells.Item(6, 10).Value = 241
$ws.Cells.Item(7, 10).Value = 514

$ws = $wb.Worksheets.Item("Englewood")
$ws.Cells.Item(2, 10).Value = 429
$ws.Cells.Item(3, 10).Value = 504
$ws.Cells.Item(6, 10).Value = 364
$ws.Cells.Item(7, 10).Value = 1426

$ws = $wb.Worksheets.Item("Chatham")
$ws.Cells.Item(2, 10).Value = 187
$ws.Cells.Item(6, 10).Value = 297
$ws.Cells.Item(7, 10).Value = 770

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Cells.Item(6, 10).Value = 58
$ws.Cells.Item(7, 10).Value = 143

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Cells.Item(3, 10).Value = 47
$ws.Cells.Item(7, 10).Value = 200

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Cells.Item(6, 10).Value = 611
$ws.Cells.Item(7, 10).Value = 1147

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Cells.Item(2, 10).Value = 79
$ws.Cells.Item(7, 10).Value = 419

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Cells.Item(3, 10).Value = 22
$ws.Cells.Item(6, 10).Value = 37
$ws.Cells.Item(7, 10).Value = 90

$ws = $wb.Worksheets.Item("Douglas")
$ws.Cells.Item(6, 10).Value = 65
$ws.Cells.Item(7, 10).Value = 243

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Cells.Item(4, 10).Value = 18
$ws.Cells.Item(6, 10).Value = 104
$ws.Cells.Item(7, 10).Value = 288

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Cells.Item(2, 10).Value = 80
$ws.Cells.Item(3, 10).Value = 126
$ws.Cells.Item(7, 10).Value = 306

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Cells.Item(2, 10).Value = 49
$ws.Cells.Item(7, 10).Value = 175

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Cells.Item(2, 10).Value = 158
$ws.Cells.Item(3, 10).Value = 189
$ws.Cells.Item(6, 10).Value = 160
$ws.Cells.Item(7, 10).Value = 562

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Cells.Item(6, 10).Value = 109
$ws.Cells.Item(7, 10).Value = 361

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Cells.Item(6, 10).Value = 47
$ws.Cells.Item(7, 10).Value = 120

$ws = $wb.Worksheets.Item("West Loop")
$ws.Cells.Item(3, 10).Value = 56
$ws.Cells.Item(7, 10).Value = 290

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Cells.Item(6, 10).Value = 145
$ws.Cells.Item(7, 10).Value = 325

$ws = $wb.Worksheets.Item("Greektown")
$ws.Cells.Item(5, 10).Value = 7
$ws.Cells.Item(6, 10).Value = 19

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Cells.Item(2, 10).Value = 132
$ws.Cells.Item(7, 10).Value = 468

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Cells.Item(2, 10).Value = 81
$ws.Cells.Item(6, 10).Value = 93

$ws = $wb.Worksheets.Item("Uptown")
$ws.Cells.Item(6, 10).Value = 101
$ws.Cells.Item(7, 10).Value = 333

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Cells.Item(6, 10).Value = 40
$ws.Cells.Item(7, 10).Value = 80

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Cells.Item(6, 10).Value = 29
$ws.Cells.Item(7, 10).Value = 168

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Cells.Item(6, 10).Value = 83
$ws.Cells.Item(7, 10).Value = 279

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Cells.Item(6, 10).Value = 133
$ws.Cells.Item(7, 10).Value = 326

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Cells.Item(2, 10).Value = 54
$ws.Cells.Item(7, 10).Value = 156

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Cells.Item(4, 10).Value = 21
$ws.Cells.Item(7, 10).Value = 227

$ws = $wb.Worksheets.Item("Clearing")
$ws.Cells.Item(6, 10).Value = 14
$ws.Cells.Item(7, 10).Value = 63

$ws = $wb.Worksheets.Item("Oakland")
$ws.Cells.Item(6, 10).Value = 35
$ws.Cells.Item(7, 10).Value = 86

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Cells.Item(6, 10).Value = 23
$ws.Cells.Item(7, 10).Value = 45

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Cells.Item(3, 10).Value = 11
$ws.Cells.Item(7, 10).Value = 86

$ws = $wb.Worksheets.Item("Millenium Park")
$ws.Cells.Item(6, 10).Value = 9
$ws.Cells.Item(7, 10).Value = 15
